{"js": "// Replace the date line and every two-digit multiplication expression\n// with the updated values from the target revision.\nconst replacements = [\n  [\"2024-04-03 Wednesday\", \"2024-04-04 Thursday\"],\n  [\"37\u00d795=\", \"46\u00d721=\"],\n  [\"73\u00d768=\", \"36\u00d786=\"],\n  [\"68\u00d743=\", \"94\u00d747=\"],\n  [\"17\u00d775=\", \"40\u00d786=\"],\n  [\"77\u00d774=\", \"61\u00d766=\"],\n  [\"15\u00d789=\", \"95\u00d742=\"],\n  [\"48\u00d746=\", \"17\u00d759=\"],\n  [\"27\u00d787=\", \"52\u00d741=\"],\n  [\"75\u00d711=\", \"22\u00d795=\"],\n  [\"66\u00d760=\", \"99\u00d743=\"],\n  [\"91\u00d798=\", \"74\u00d719=\"],\n  [\"82\u00d740=\", \"22\u00d748=\"],\n  [\"56\u00d780=\", \"60\u00d749=\"],\n  [\"81\u00d784=\", \"83\u00d764=\"],\n  [\"99\u00d735=\", \"51\u00d733=\"],\n  [\"79\u00d742=\", \"44\u00d738=\"],\n  [\"72\u00d780=\", \"32\u00d781=\"],\n  [\"28\u00d745=\", \"13\u00d789=\"],\n  [\"53\u00d760=\", \"97\u00d779=\"],\n  [\"62\u00d737=\", \"12\u00d725=\"],\n  [\"24\u00d780=\", \"35\u00d784=\"],\n  [\"57\u00d797=\", \"11\u00d747=\"],\n  [\"74\u00d752=\", \"18\u00d798=\"],\n  [\"17\u00d730=\", \"75\u00d715=\"],\n  [\"82\u00d743=\", \"95\u00d766=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit multiplication expression\n# to match the target revision.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-04-03 Wednesday\", \"2024-04-04 Thursday\"),\n  @(\"37\u00d795=\", \"46\u00d721=\"),\n  @(\"73\u00d768=\", \"36\u00d786=\"),\n  @(\"68\u00d743=\", \"94\u00d747=\"),\n  @(\"17\u00d775=\", \"40\u00d786=\"),\n  @(\"77\u00d774=\", \"61\u00d766=\"),\n  @(\"15\u00d789=\", \"95\u00d742=\"),\n  @(\"48\u00d746=\", \"17\u00d759=\"),\n  @(\"27\u00d787=\", \"52\u00d741=\"),\n  @(\"75\u00d711=\", \"22\u00d795=\"),\n  @(\"66\u00d760=\", \"99\u00d743=\"),\n  @(\"91\u00d798=\", \"74\u00d719=\"),\n  @(\"82\u00d740=\", \"22\u00d748=\"),\n  @(\"56\u00d780=\", \"60\u00d749=\"),\n  @(\"81\u00d784=\", \"83\u00d764=\"),\n  @(\"99\u00d735=\", \"51\u00d733=\"),\n  @(\"79\u00d742=\", \"44\u00d738=\"),\n  @(\"72\u00d780=\", \"32\u00d781=\"),\n  @(\"28\u00d745=\", \"13\u00d789=\"),\n  @(\"53\u00d760=\", \"97\u00d779=\"),\n  @(\"62\u00d737=\", \"12\u00d725=\"),\n  @(\"24\u00d780=\", \"35\u00d784=\"),\n  @(\"57\u00d797=\", \"11\u00d747=\"),\n  @(\"74\u00d752=\", \"18\u00d798=\"),\n  @(\"17\u00d730=\", \"75\u00d715=\"),\n  @(\"82\u00d743=\", \"95\u00d766=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
